$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57: AF57 gains a fractional remainder, and the trailing AZ:BT block
# (a run of zeros) is cleared out entirely.
$ws.Range("AF57").Value = 44888.008
$ws.Range("AZ57:BT57").ClearContents()

# Row 58: trailing AZ:BT zero block cleared.
$ws.Range("AZ58:BT58").ClearContents()

# Row 64: AE:AZ zero block cleared (BA64 onward stays as-is).
$ws.Range("AE64:AZ64").ClearContents()

# Row 71: trailing AZ:BT zero block cleared.
$ws.Range("AZ71:BT71").ClearContents()

# Row 72: trailing AZ:BT zero block cleared.
$ws.Range("AZ72:BT72").ClearContents()

# Row 73: trailing AZ:BT zero block cleared.
$ws.Range("AZ73:BT73").ClearContents()

# Row 77: trailing AZ:BT zero block cleared.
$ws.Range("AZ77:BT77").ClearContents()

# Row 78: trailing AZ:BT zero block cleared.
$ws.Range("AZ78:BT78").ClearContents()

# Row 79: AE:BT zero block cleared entirely.
$ws.Range("AE79:BT79").ClearContents()
